$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54: Added 3rd Sonia Delaunay - 23/09/2019 - 0.5 hrs @ 25
$ws.Range("A54").Value = "Added 3rd Sonia Delaunay"
$ws.Range("B54").Value = 43731
$ws.Range("C54").Value = 0.5
$ws.Range("D54").Value = 25

# Row 55: Added Miro Prints - 23/09/2019 - 2 hrs @ 25
$ws.Range("A55").Value = "Added Miro Prints"
$ws.Range("B55").Value = 43731
$ws.Range("C55").Value = 2
$ws.Range("D55").Value = 25

# Restore the selection/view state left by the author after the edit
$ws.Range("D56").Select()
